# Adds the "Environment" / "Watershed" sub-module rows (28-33) describing
# the computation of watershed concentration, together with the new
# shared strings, styles and data-validation ranges that go with them.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EstimatedVariables")

# ---------------------------------------------------------------------
# New row data. Each entry: row number, Sub-module, Variable, Description,
# Units, Distribution, Truncated, Point-estimate value
# ---------------------------------------------------------------------

$rows = @(
    @{ Row = 28; Variable = "m_apply";  Description = "Manure application rate ";              Units = "kg/m²"; Value = 10 },
    @{ Row = 29; Variable = "runoff";   Description = "Fraction of manure that runs off";       Units = "decimal"; Value = 0.1 },
    @{ Row = 30; Variable = "transport";Description = "Fraction of E.coli surviving transport"; Units = "decimal"; Value = 0.5 },
    @{ Row = 31; Variable = "w_area";   Description = "Surface area of the waterbody";          Units = "m2"; Value = 10000 },
    @{ Row = 32; Variable = "depth";    Description = "Average depth of the waterbody";         Units = "m"; Value = 2 },
    @{ Row = 33; Variable = "c_factor"; Description = "Conversion factor from kg to g";         Units = "decimal"; Value = 1000000 }
)

foreach ($r in $rows) {
    $i = $r.Row

    $ws.Cells.Item($i, 1).Value = "Environment"
    $ws.Cells.Item($i, 2).Value = "Watershed"
    $ws.Cells.Item($i, 3).Value = $r.Variable
    $ws.Cells.Item($i, 4).Value = $r.Description
    $ws.Cells.Item($i, 5).Value = $r.Units
    $ws.Cells.Item($i, 6).Value = "Point estimate"
    $ws.Cells.Item($i, 7).Value = "No"
    $ws.Cells.Item($i, 12).Value = $r.Value
    $ws.Cells.Item($i, 15).Value = "User input"
}

# ---------------------------------------------------------------------
# Styling to match the rest of the "Sub-module" column blocks:
#  - Column A ("Environment"): bold-ish header style with a themed fill
#  - Column B ("Watershed"): plain cell with a lighter themed fill
#  - Column D (description): italic font
# ---------------------------------------------------------------------

$ws.Range("A28:A33").Font.Italic = $false
$ws.Range("A28:A33").Interior.ThemeColor = 8
$ws.Range("B28:B33").Interior.ThemeColor = 8
$ws.Range("D28:D33").Font.Italic = $true

# Truncated ("No") column keeps the same style used elsewhere (G column)
$ws.Range("G28:G33").Style = $ws.Range("G27").Style

# ---------------------------------------------------------------------
# Grow the worksheet dimension / selection to include the new rows
# ---------------------------------------------------------------------

$ws.Range("D33").Select()

# ---------------------------------------------------------------------
# Extend the data validation ranges so the new rows get the same
# dropdowns as the rest of the table
# ---------------------------------------------------------------------

$ws.Range("G2:G33").Validation.Delete()
$ws.Range("G2:G33").Validation.Add(3, 1, 1, "=ListItem!`$B`$3:`$B`$4")

$ws.Range("F6:F7").Validation.Delete()
$ws.Range("F6:F7").Validation.Add(3, 1, 1, "=ListItem!`$A`$2:`$A`$17")

$ws.Range("F2:F5").Validation.Delete()
$ws.Range("F2:F5").Validation.Add(3, 1, 1, "=ListItem!`$A`$2:`$A`$17")

Write-Host "done"
